# Rename the right-shift (RSHT) and left-shift (LSHT) mnemonics in the
# instruction-set reference table to their new short forms SHR / SHL.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "SHR"
$ws.Range("A18").Value = "SHL"

# Move the active selection, matching the author's final cursor position.
$ws.Range("F11").Select() | Out-Null
